$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column N values for rows 11-25 to reflect revised figures
$ws.Range("N11").Value = 240
$ws.Range("N12").Value = 240
$ws.Range("N13").Value = 300
$ws.Range("N14").Value = 2000
$ws.Range("N15").Value = 2000
$ws.Range("N16").Value = 120
$ws.Range("N17").Value = 120
$ws.Range("N18").Value = 240
$ws.Range("N19").Value = 200
$ws.Range("N20").Value = 1200
$ws.Range("N21").Value = 2000
$ws.Range("N22").Value = 2000
$ws.Range("N23").Value = 2000
$ws.Range("N24").Value = 2000
$ws.Range("N25").Value = 2000

# Move the active cell selection from T11 to T13, matching the saved view state
$ws.Range("T13").Select()
